# Data update using git
# Apply numeric corrections to the "Inscritos", "Pagos" and
# "Inscrições homologadas" columns on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value  = 30
$ws.Range("F6").Value  = 27
$ws.Range("H6").Value  = 34
$ws.Range("E10").Value = 32
$ws.Range("F11").Value = 14
$ws.Range("H11").Value = 15
$ws.Range("E15").Value = 103
$ws.Range("E17").Value = 27
